$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 300
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -16
$ws.Range("N18").Value = ""

$ws.Range("H41").Value = 466.66666
$ws.Range("I41").Value = 466.66666
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 466.66666
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -26.66665999999998

$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968

$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 4000
$ws.Range("N64").Value = -4496
$ws.Range("M64").Value = ""

$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 4000
$ws.Range("N67").Value = -5716
$ws.Range("M67").Value = ""

$ws.Range("H70").Value = 3125
$ws.Range("I70").Value = 4150
$ws.Range("J70").Value = 2100
$ws.Range("K70").Value = 12450
$ws.Range("L70").Value = 6300
$ws.Range("M70").Value = -12180
$ws.Range("N70").Value = -6840

$ws.Range("H73").Value = 3125
$ws.Range("I73").Value = 4150
$ws.Range("J73").Value = 2100
$ws.Range("K73").Value = 12450
$ws.Range("L73").Value = 6300
$ws.Range("M73").Value = -11514
$ws.Range("N73").Value = -8172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1944.2222
$ws.Range("I2").Value = 2062.375
$ws.Range("J2").Value = 999
$ws.Range("K2").Value = 2062.375
$ws.Range("L2").Value = 999
$ws.Range("M2").Value = -1949.375
$ws.Range("N2").Value = -1225

$ws.Range("H45").Value = 4645.6665
$ws.Range("I45").Value = 1624.8334
$ws.Range("J45").Value = 7666.5
$ws.Range("K45").Value = 1624.8334
$ws.Range("L45").Value = 7666.5
$ws.Range("M45").Value = -1247.8334
$ws.Range("N45").Value = -8420.5

$ws.Range("H63").Value = 6005.1113
$ws.Range("I63").Value = 5409.8
$ws.Range("J63").Value = 6749.25
$ws.Range("K63").Value = 5409.8
$ws.Range("L63").Value = 6749.25
$ws.Range("M63").Value = -4723.8
$ws.Range("N63").Value = -8121.25

$ws.Range("H66").Value = 6005.1113
$ws.Range("I66").Value = 5409.8
$ws.Range("J66").Value = 6749.25
$ws.Range("K66").Value = 27049
$ws.Range("L66").Value = 33746.25
$ws.Range("M66").Value = -23617
$ws.Range("N66").Value = -40610.25

$ws.Range("H88").Value = 1084
$ws.Range("I88").Value = 800
$ws.Range("J88").Value = 1226
$ws.Range("K88").Value = 800
$ws.Range("L88").Value = 1226
$ws.Range("M88").Value = -394
$ws.Range("N88").Value = -2038

$ws.Range("H91").Value = 1084
$ws.Range("I91").Value = 800
$ws.Range("J91").Value = 1226
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 1226
$ws.Range("M91").Value = 604
$ws.Range("N91").Value = -4034

$ws.Range("H98").Value = 31782.715
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 31782.715
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 31782.715
$ws.Range("N98").Value = -37772.715

$ws.Range("H116").Value = 1944.2222
$ws.Range("I116").Value = 2062.375
$ws.Range("J116").Value = 999
$ws.Range("K116").Value = 2062.375
$ws.Range("L116").Value = 999
$ws.Range("M116").Value = 231.625
$ws.Range("N116").Value = -5587

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1944.2222
$ws.Range("I3").Value = 2062.375
$ws.Range("J3").Value = 999
$ws.Range("K3").Value = 2062.375
$ws.Range("L3").Value = 999
$ws.Range("M3").Value = -1948.375
$ws.Range("N3").Value = -1227

$ws.Range("H54").Value = 6000
$ws.Range("I54").Value = 6000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 6000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -5516

$ws.Range("H86").Value = 809.5
$ws.Range("I86").Value = 750
$ws.Range("J86").Value = 898.75
$ws.Range("K86").Value = 750
$ws.Range("L86").Value = 898.75
$ws.Range("M86").Value = 373
$ws.Range("N86").Value = -3144.75

$ws.Range("H89").Value = 809.5
$ws.Range("I89").Value = 750
$ws.Range("J89").Value = 898.75
$ws.Range("K89").Value = 3750
$ws.Range("L89").Value = 4493.75
$ws.Range("M89").Value = 1866
$ws.Range("N89").Value = -15725.75

$ws.Range("H94").Value = 838.46155
$ws.Range("I94").Value = 824.25
$ws.Range("J94").Value = 1009
$ws.Range("K94").Value = 824.25
$ws.Range("L94").Value = 1009
$ws.Range("M94").Value = -373.25
$ws.Range("N94").Value = -1911

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 309.6
$ws.Range("I22").Value = 284.5
$ws.Range("J22").Value = 359.8
$ws.Range("K22").Value = 284.5
$ws.Range("L22").Value = 359.8
$ws.Range("M22").Value = 65.5
$ws.Range("N22").Value = -1059.8

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = ""

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = ""

$ws.Range("H99").Value = 1200
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 298

$ws.Range("H122").Value = 1311.125
$ws.Range("I122").Value = 1311.125
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3933.375
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1483.375

$ws.Range("H126").Value = 1200
$ws.Range("I126").Value = 1200
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3600
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -888

$ws.Range("H43").Value = 10315.667
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10315.667
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 10315.667
$ws.Range("N43").Value = -10617.667

$ws.Range("H80").Value = 2167
$ws.Range("I80").Value = 2329
$ws.Range("J80").Value = 2005
$ws.Range("K80").Value = 2329
$ws.Range("L80").Value = 2005
$ws.Range("M80").Value = -1331
$ws.Range("N80").Value = -4001

$ws.Range("H83").Value = 2167
$ws.Range("I83").Value = 2329
$ws.Range("J83").Value = 2005
$ws.Range("K83").Value = 11645
$ws.Range("L83").Value = 10025
$ws.Range("M83").Value = -6653
$ws.Range("N83").Value = -20009

$ws.Range("H123").Value = 74354.164
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 74354.164
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 74354.164
$ws.Range("N123").Value = -79254.164

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1479.6
$ws.Range("I16").Value = 1479.6
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1479.6
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1309.6

$ws.Range("H22").Value = 1843.5
$ws.Range("I22").Value = 487.5
$ws.Range("J22").Value = 3199.5
$ws.Range("K22").Value = 487.5
$ws.Range("L22").Value = 3199.5
$ws.Range("M22").Value = -192.5
$ws.Range("N22").Value = -3789.5

$ws.Range("H27").Value = 1843.5
$ws.Range("I27").Value = 487.5
$ws.Range("J27").Value = 3199.5
$ws.Range("K27").Value = 487.5
$ws.Range("L27").Value = 3199.5
$ws.Range("M27").Value = -380.5
$ws.Range("N27").Value = -3413.5

$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -3376

$ws.Range("H68").Value = 3108.3635
$ws.Range("I68").Value = 2849.375
$ws.Range("J68").Value = 3799
$ws.Range("K68").Value = 2849.375
$ws.Range("L68").Value = 3799
$ws.Range("M68").Value = -2100.375
$ws.Range("N68").Value = -5297

$ws.Range("H71").Value = 3108.3635
$ws.Range("I71").Value = 2849.375
$ws.Range("J71").Value = 3799
$ws.Range("K71").Value = 14246.875
$ws.Range("L71").Value = 18995
$ws.Range("M71").Value = -10502.875
$ws.Range("N71").Value = -26483

$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12000
$ws.Range("N122").Value = -16900
$ws.Range("M122").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2000000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 2000000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 2000000
$ws.Range("N3").Value = -2000228

$ws.Range("H14").Value = 1699
$ws.Range("I14").Value = 1165
$ws.Range("J14").Value = 2500
$ws.Range("K14").Value = 1165
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = -997
$ws.Range("N14").Value = -2836

$ws.Range("H41").Value = 49654.668
$ws.Range("I41").Value = 46590.6
$ws.Range("J41").Value = 64975
$ws.Range("K41").Value = 46590.6
$ws.Range("L41").Value = 64975
$ws.Range("M41").Value = -46200.6
$ws.Range("N41").Value = -65755

$ws.Range("H113").Value = 1719.5714
$ws.Range("I113").Value = 1172.8334
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 3518.5002
$ws.Range("L113").Value = 15000
$ws.Range("M113").Value = -1348.5002
$ws.Range("N113").Value = -19340

$ws.Range("H132").Value = 5443.5
$ws.Range("I132").Value = 2258
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 6774
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -4244
$ws.Range("N132").Value = -50060
